$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-09-06"

# Update the header label cell (I1) that says "2022 (through 09-05)"
$ws.Range("I1").Value = "2022 (through 09-06)"

# Update the updated/corrected data values
$ws.Range("I9").Value = 166
$ws.Range("I10").Value = 31
$ws.Range("I14").Value = 1168
